$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 59.857143
$ws.Range("I5").Value = 53.166668
$ws.Range("K5").Value = 53.166668
$ws.Range("M5").Value = 61.833332
$ws.Range("H53").Value = 1297.3529
$ws.Range("I53").Value = 132.33333
$ws.Range("J53").Value = 1932.8182
$ws.Range("K53").Value = 132.33333
$ws.Range("L53").Value = 1932.8182
$ws.Range("M53").Value = 504.66667
$ws.Range("N53").Value = -3206.8182
$ws.Range("H112").Value = 1889.4286
$ws.Range("J112").Value = 2457.25
$ws.Range("L112").Value = 7371.75
$ws.Range("N112").Value = -9587.75
$ws.Range("H133").Value = 60172.75
$ws.Range("J133").Value = 60172.75
$ws.Range("L133").Value = 60172.75
$ws.Range("N133").Value = -70292.75
$ws.Range("H138").Value = 2719
$ws.Range("I138").Value = 1655.2667
$ws.Range("J138").Value = 3581.4866
$ws.Range("K138").Value = 4965.800099999999
$ws.Range("L138").Value = 10744.4598
$ws.Range("M138").Value = 174.1999000000005
$ws.Range("N138").Value = -21024.4598

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2436.4092
$ws.Range("I32").Value = 1704.225
$ws.Range("K32").Value = 1704.225
$ws.Range("M32").Value = -1417.225
$ws.Range("H55").Value = 6512
$ws.Range("I55").Value = 6512
$ws.Range("K55").Value = 6512
$ws.Range("M55").Value = -6197
$ws.Range("H61").Value = 7113.15
$ws.Range("I61").Value = 4017.6667
$ws.Range("J61").Value = 16399.6
$ws.Range("K61").Value = 4017.6667
$ws.Range("L61").Value = 16399.6
$ws.Range("M61").Value = -3805.6667
$ws.Range("N61").Value = -16823.6
$ws.Range("H63").Value = 3999.5
$ws.Range("I63").Value = 3999.5
$ws.Range("K63").Value = 3999.5
$ws.Range("M63").Value = -3313.5
$ws.Range("H66").Value = 3999.5
$ws.Range("I66").Value = 3999.5
$ws.Range("K66").Value = 19997.5
$ws.Range("M66").Value = -16565.5
$ws.Range("H74").Value = 27780142
$ws.Range("I74").Value = 33335272
$ws.Range("K74").Value = 33335272
$ws.Range("M74").Value = -33334398
$ws.Range("H77").Value = 27780142
$ws.Range("I77").Value = 33335272
$ws.Range("K77").Value = 166676360
$ws.Range("M77").Value = -166671992
$ws.Range("H122").Value = 3337.1667
$ws.Range("I122").Value = 2731.4546
$ws.Range("K122").Value = 8194.363799999999
$ws.Range("M122").Value = -5744.363799999999
$ws.Range("H132").Value = 4197.891
$ws.Range("I132").Value = 2898.25
$ws.Range("K132").Value = 8694.75
$ws.Range("M132").Value = -6164.75
$ws.Range("H136").Value = 7113.15
$ws.Range("I136").Value = 4017.6667
$ws.Range("J136").Value = 16399.6
$ws.Range("K136").Value = 12053.0001
$ws.Range("L136").Value = 49198.8
$ws.Range("M136").Value = -9503.000100000001
$ws.Range("N136").Value = -54298.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3569.8667
$ws.Range("I20").Value = 2481.8333
$ws.Range("K20").Value = 2481.8333
$ws.Range("M20").Value = -2234.8333
$ws.Range("H81").Value = 61822.4
$ws.Range("J81").Value = 61822.4
$ws.Range("L81").Value = 61822.4
$ws.Range("N81").Value = -63944.4
$ws.Range("H84").Value = 61822.4
$ws.Range("J84").Value = 61822.4
$ws.Range("L84").Value = 185467.2
$ws.Range("N84").Value = -196075.2
$ws.Range("H105").Value = 16670.055
$ws.Range("I105").Value = 17076.215
$ws.Range("K105").Value = 17076.215
$ws.Range("M105").Value = -15329.215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 124995
$ws.Range("J9").Value = 124995
$ws.Range("L9").Value = 124995
$ws.Range("N9").Value = -125331
$ws.Range("H99").Value = 3456
$ws.Range("I99").Value = 2912
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 2912
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -1414
$ws.Range("N99").Value = -6996
$ws.Range("H105").Value = 4425.75
$ws.Range("I105").Value = 2342.1428
$ws.Range("K105").Value = 2342.1428
$ws.Range("M105").Value = -595.1428000000001
$ws.Range("H112").Value = 67112.8
$ws.Range("J112").Value = 67112.8
$ws.Range("L112").Value = 67112.8
$ws.Range("N112").Value = -70066.8
$ws.Range("H126").Value = 3456
$ws.Range("I126").Value = 2912
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 8736
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -6266
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 10354.363
$ws.Range("J32").Value = 10354.363
$ws.Range("L32").Value = 31063.089
$ws.Range("N32").Value = -31629.089
$ws.Range("H38").Value = 17.5
$ws.Range("J38").Value = 15
$ws.Range("L38").Value = 45
$ws.Range("N38").Value = -739
$ws.Range("H42").Value = 7849
$ws.Range("J42").Value = 7849
$ws.Range("L42").Value = 23547
$ws.Range("N42").Value = -24615
$ws.Range("H68").Value = 3592.652
$ws.Range("I68").Value = 1400
$ws.Range("J68").Value = 3692.318
$ws.Range("K68").Value = 4200
$ws.Range("L68").Value = 11076.954
$ws.Range("M68").Value = -3389
$ws.Range("N68").Value = -12698.954
$ws.Range("H71").Value = 3592.652
$ws.Range("I71").Value = 1400
$ws.Range("J71").Value = 3692.318
$ws.Range("K71").Value = 12600
$ws.Range("L71").Value = 33230.862
$ws.Range("M71").Value = -8544
$ws.Range("N71").Value = -41342.862
$ws.Range("H103").Value = 1895.375
$ws.Range("J103").Value = 2165.6
$ws.Range("L103").Value = 6496.799999999999
$ws.Range("N103").Value = -8254.799999999999
$ws.Range("H132").Value = 4692.4116
$ws.Range("I132").Value = 2921.2222
$ws.Range("J132").Value = 6685
$ws.Range("K132").Value = 26290.9998
$ws.Range("L132").Value = 60165
$ws.Range("M132").Value = -23760.9998
$ws.Range("N132").Value = -65225

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 35500
$ws.Range("J104").Value = 35500
$ws.Range("L104").Value = 35500
$ws.Range("N104").Value = -42488
$ws.Range("H126").Value = 12671.333
$ws.Range("J126").Value = 12671.333
$ws.Range("L126").Value = 38013.999
$ws.Range("N126").Value = -42953.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 20865.666
$ws.Range("I7").Value = 3297.8
$ws.Range("K7").Value = 3297.8
$ws.Range("M7").Value = -3185.8
$ws.Range("H40").Value = 13314.728
$ws.Range("J40").Value = 18335
$ws.Range("L40").Value = 18335
$ws.Range("N40").Value = -18607
$ws.Range("H93").Value = 1292.826
$ws.Range("I93").Value = 1336.4286
$ws.Range("K93").Value = 1336.4286
$ws.Range("M93").Value = -88.42859999999996
$ws.Range("H126").Value = 20865.666
$ws.Range("I126").Value = 3297.8
$ws.Range("K126").Value = 9893.400000000001
$ws.Range("M126").Value = -7423.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7805.0415
$ws.Range("I122").Value = 1887.7333
$ws.Range("K122").Value = 5663.199900000001
$ws.Range("M122").Value = -3213.199900000001
$ws.Range("H126").Value = 3923.7
$ws.Range("I126").Value = 3277.8823
$ws.Range("J126").Value = 7583.3335
$ws.Range("K126").Value = 9833.6469
$ws.Range("L126").Value = 22750.0005
$ws.Range("M126").Value = -7363.6469
$ws.Range("N126").Value = -27690.0005
$ws.Range("H135").Value = 53019.1
$ws.Range("J135").Value = 53019.1
$ws.Range("L135").Value = 53019.1
$ws.Range("N135").Value = -63159.1
$ws.Range("H136").Value = 3811.0205
$ws.Range("I136").Value = 2646.8096
$ws.Range("J136").Value = 10796.286
$ws.Range("K136").Value = 7940.4288
$ws.Range("L136").Value = 32388.858
$ws.Range("M136").Value = -5390.4288
$ws.Range("N136").Value = -37488.858
